$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F2 value from "Stor_UEF" to "Instant_UEF" (renaming cohort value)
$ws.Range("F2").Value = "Instant_UEF"

# Update the current selection to F2 to match the saved selection state
$ws.Range("F2").Select()
